# Create the new "PO Forecast" worksheet at the end of the workbook
$wb = $excel.ActiveWorkbook

# --- Update header labels on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Copy the header style (bold, bordered, centered) from an existing header cell
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Copy the date-column style from an existing date cell
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A11").PasteSpecial(-4122)

# --- Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- Data rows ---
$wsForecast.Range("A2").Value = 45592.99999999999
$wsForecast.Range("B2").Value = 40
$wsForecast.Range("C2").Value = 40.00008001630086
$wsForecast.Range("D2").Value = 40.00008001730938
$wsForecast.Range("A3").Value = 45599.99999999999
$wsForecast.Range("B3").Value = 130
$wsForecast.Range("C3").Value = 130.0001199957021
$wsForecast.Range("D3").Value = 130.0001199966558
$wsForecast.Range("A4").Value = 45606.99999999999
$wsForecast.Range("B4").Value = 220
$wsForecast.Range("C4").Value = 220.0001589354084
$wsForecast.Range("D4").Value = 220.0001609765814
$wsForecast.Range("A5").Value = 45613.99999999999
$wsForecast.Range("B5").Value = 310
$wsForecast.Range("C5").Value = 310.0001965682857
$wsForecast.Range("D5").Value = 310.0002029249627
$wsForecast.Range("A6").Value = 45620.99999999999
$wsForecast.Range("B6").Value = 400
$wsForecast.Range("C6").Value = 400.000233285365
$wsForecast.Range("D6").Value = 400.0002457477284
$wsForecast.Range("A7").Value = 45627.99999999999
$wsForecast.Range("B7").Value = 490
$wsForecast.Range("C7").Value = 490.0002692279894
$wsForecast.Range("D7").Value = 490.0002892726195
$wsForecast.Range("A8").Value = 45634.99999999999
$wsForecast.Range("B8").Value = 580
$wsForecast.Range("C8").Value = 580.0003044742665
$wsForecast.Range("D8").Value = 580.0003339675025
$wsForecast.Range("A9").Value = 45641.99999999999
$wsForecast.Range("B9").Value = 670
$wsForecast.Range("C9").Value = 670.0003395008495
$wsForecast.Range("D9").Value = 670.0003785859716
$wsForecast.Range("A10").Value = 45648.99999999999
$wsForecast.Range("B10").Value = 760
$wsForecast.Range("C10").Value = 760.0003741988096
$wsForecast.Range("D10").Value = 760.0004231724661
$wsForecast.Range("A11").Value = 45655.99999999999
$wsForecast.Range("B11").Value = 850
$wsForecast.Range("C11").Value = 850.0004085020474
$wsForecast.Range("D11").Value = 850.0004684128229
